$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 3045.8
$ws.Range("I10").Value = 3502
$ws.Range("J10").Value = 2741.6667
$ws.Range("K10").Value = 3502
$ws.Range("L10").Value = 2741.6667
$ws.Range("M10").Value = -3209
$ws.Range("N10").Value = -3327.6667

$ws.Range("H86").Value = 5791.857
$ws.Range("I86").Value = 6183
$ws.Range("K86").Value = 6183
$ws.Range("M86").Value = -5060

$ws.Range("H89").Value = 5791.857
$ws.Range("I89").Value = 6183
$ws.Range("K89").Value = 30915
$ws.Range("M89").Value = -25299

$ws.Range("H112").Value = 2776.3
$ws.Range("J112").Value = 3502
$ws.Range("L112").Value = 10506
$ws.Range("N112").Value = -12722

$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -1746

$ws.Range("H125").Value = 20151
$ws.Range("I125").Value = 877
$ws.Range("J125").Value = 26575.666
$ws.Range("K125").Value = 7893
$ws.Range("L125").Value = 239180.994
$ws.Range("M125").Value = -5433
$ws.Range("N125").Value = -244100.994

$ws.Range("H131").Value = 46999.4
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 46999.4
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 140998.2
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -151078.2

$ws.Range("H132").Value = 6637.3022
$ws.Range("I132").Value = 6530.8975
$ws.Range("K132").Value = 19592.6925
$ws.Range("M132").Value = -17062.6925

$ws.Range("H137").Value = 16133547
$ws.Range("I137").Value = 37038708
$ws.Range("J137").Value = 6708.086
$ws.Range("K137").Value = 111116124
$ws.Range("L137").Value = 20124.258
$ws.Range("M137").Value = -111113574
$ws.Range("N137").Value = -25224.258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 12999.667
$ws.Range("I28").Value = 12999.667
$ws.Range("K28").Value = 12999.667
$ws.Range("M28").Value = -12807.667

$ws.Range("H32").Value = 139134.33
$ws.Range("I32").Value = 209152.67
$ws.Range("K32").Value = 209152.67
$ws.Range("M32").Value = -208865.67

$ws.Range("H61").Value = 2568308.5
$ws.Range("I61").Value = 4207.2285
$ws.Range("K61").Value = 4207.2285
$ws.Range("M61").Value = -3995.2285

$ws.Range("H97").Value = 125003816
$ws.Range("I97").Value = 3661.4
$ws.Range("K97").Value = 3661.4
$ws.Range("M97").Value = -3165.4

$ws.Range("H99").Value = 12999.667
$ws.Range("I99").Value = 12999.667
$ws.Range("K99").Value = 12999.667
$ws.Range("M99").Value = -10004.667

$ws.Range("H136").Value = 2568308.5
$ws.Range("I136").Value = 4207.2285
$ws.Range("K136").Value = 12621.6855
$ws.Range("M136").Value = -10071.6855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 13978.75
$ws.Range("J99").Value = 1799.6
$ws.Range("L99").Value = 1799.6
$ws.Range("N99").Value = -4795.6

$ws.Range("H134").Value = 2320176.8
$ws.Range("I134").Value = 3603.7544
$ws.Range("J134").Value = 11123154
$ws.Range("K134").Value = 10811.2632
$ws.Range("L134").Value = 33369462
$ws.Range("M134").Value = -8276.263199999999
$ws.Range("N134").Value = -33374532

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H94").Value = 9254.786
$ws.Range("I94").Value = 21797
$ws.Range("J94").Value = 2286.889
$ws.Range("K94").Value = 21797
$ws.Range("L94").Value = 2286.889
$ws.Range("M94").Value = -21346
$ws.Range("N94").Value = -3188.889

$ws.Range("H107").Value = 667.5
$ws.Range("I107").Value = 554.5
$ws.Range("K107").Value = 554.5
$ws.Range("M107").Value = 1365.5

$ws.Range("H140").Value = 99993
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 125.521736
$ws.Range("J2").Value = 161.63637
$ws.Range("L2").Value = 969.81822
$ws.Range("N2").Value = -1195.81822

$ws.Range("H61").Value = 335.83334
$ws.Range("J61").Value = 333
$ws.Range("L61").Value = 999
$ws.Range("N61").Value = -1429

$ws.Range("H98").Value = 361.25
$ws.Range("I98").Value = 446.33334
$ws.Range("J98").Value = 310.2
$ws.Range("K98").Value = 1339.00002
$ws.Range("L98").Value = 930.5999999999999
$ws.Range("M98").Value = 158.9999800000001
$ws.Range("N98").Value = -3926.6

$ws.Range("H139").Value = 6180
$ws.Range("I139").Value = 2700.077
$ws.Range("K139").Value = 8100.231000000001
$ws.Range("M139").Value = -2960.231000000001

$ws.Range("H140").Value = 10225
$ws.Range("I140").Value = 9000
$ws.Range("J140").Value = 11450
$ws.Range("K140").Value = 27000
$ws.Range("L140").Value = 34350
$ws.Range("M140").Value = -21820
$ws.Range("N140").Value = -44710

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2802.6667
$ws.Range("J80").Value = 2996
$ws.Range("L80").Value = 2996
$ws.Range("N80").Value = -4992

$ws.Range("H83").Value = 2802.6667
$ws.Range("J83").Value = 2996
$ws.Range("L83").Value = 14980
$ws.Range("N83").Value = -24964

$ws.Range("H122").Value = 4872.8696
$ws.Range("I122").Value = 4734
$ws.Range("K122").Value = 14202
$ws.Range("M122").Value = -11752

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1892.8125
$ws.Range("J22").Value = 2806.875
$ws.Range("L22").Value = 2806.875
$ws.Range("N22").Value = -3396.875

$ws.Range("H25").Value = 2750
$ws.Range("J25").Value = 2750
$ws.Range("L25").Value = 2750
$ws.Range("N25").Value = -3210

$ws.Range("H27").Value = 1892.8125
$ws.Range("J27").Value = 2806.875
$ws.Range("L27").Value = 2806.875
$ws.Range("N27").Value = -3020.875

$ws.Range("H82").Value = 1300.4286
$ws.Range("I82").Value = 1337.24
$ws.Range("J82").Value = 993.6667
$ws.Range("K82").Value = 1337.24
$ws.Range("L82").Value = 993.6667
$ws.Range("M82").Value = -976.24
$ws.Range("N82").Value = -1715.6667

$ws.Range("H85").Value = 1300.4286
$ws.Range("I85").Value = 1337.24
$ws.Range("J85").Value = 993.6667
$ws.Range("K85").Value = 1337.24
$ws.Range("L85").Value = 993.6667
$ws.Range("M85").Value = -89.24000000000001
$ws.Range("N85").Value = -3489.6667

$ws.Range("H100").Value = 3571.2856
$ws.Range("I100").Value = 2500
$ws.Range("J100").Value = 4999.6665
$ws.Range("K100").Value = 2500
$ws.Range("L100").Value = 4999.6665
$ws.Range("M100").Value = -1959
$ws.Range("N100").Value = -6081.6665

$ws.Range("H132").Value = 4330794.5
$ws.Range("I132").Value = 6876148
$ws.Range("K132").Value = 20628444
$ws.Range("M132").Value = -20625914

$ws.Range("H136").Value = 10007817
$ws.Range("I136").Value = 11366860
$ws.Range("K136").Value = 34100580
$ws.Range("M136").Value = -34098030

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3877347.5
$ws.Range("I132").Value = 4274711
$ws.Range("J132").Value = 3050
$ws.Range("K132").Value = 12824133
$ws.Range("L132").Value = 9150
$ws.Range("M132").Value = -12821603
$ws.Range("N132").Value = -14210
